$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Y2").Value = -2.87413053845901
$ws.Range("Z2").Value = 1.484791740628893
$ws.Range("D3").Value = 18.37
$ws.Range("E3").Value = -101.63
$ws.Range("F3").Value = -1.0231506773862
$ws.Range("G3").Value = -0.4374395176060744
$ws.Range("H3").Value = -0.5059674056523835
$ws.Range("I3").Value = -0.07214517640069329
$ws.Range("J3").Value = 14.23084448760408
$ws.Range("K3").Value = 14.23084448760408
$ws.Range("L3").Value = 14.23084448760408
$ws.Range("M3").Value = -0.5851395448336352
$ws.Range("N3").Value = -0.1280656518002171
$ws.Range("O3").Value = 14.23084448760408
$ws.Range("P3").Value = 14.23084448760408
$ws.Range("Q3").Value = -0.5059674056523835
$ws.Range("R3").Value = -0.07214517640069329
$ws.Range("S3").Value = -1.0231506773862
$ws.Range("T3").Value = -0.4374395176060744
$ws.Range("U3").Value = 0.009133242001687534
$ws.Range("V3").Value = 0.02112047355486455
$ws.Range("W3").Value = 0.03439635526781767
$ws.Range("X3").Value = -0.01746109989614993
$ws.Range("Y3").Value = -0.63641847061382
$ws.Range("Z3").Value = 0.3696271059445797
$ws.Range("AA3").Value = 0.008860829757188585
$ws.Range("AB3").Value = 0.01969550581006432
$ws.Range("AC3").Value = 0.01970058862469999
$ws.Range("AD3").Value = 0.008861801443002755
$ws.Range("D4").Value = 18.53
$ws.Range("E4").Value = -81.47
$ws.Range("F4").Value = -0.8236638877349185
$ws.Range("G4").Value = -0.2921931896054355
$ws.Range("H4").Value = -0.3926778279567383
$ws.Range("I4").Value = 0.01221876139904882
$ws.Range("J4").Value = 12.6502837182674
$ws.Range("K4").Value = 12.6502837182674
$ws.Range("L4").Value = 12.6502837182674
$ws.Range("M4").Value = -0.4725395448336351
$ws.Range("N4").Value = -0.04418877312208214
$ws.Range("O4").Value = 12.6502837182674
$ws.Range("P4").Value = 12.6502837182674
$ws.Range("Q4").Value = -0.3926778279567383
$ws.Range("R4").Value = 0.01221876139904882
$ws.Range("S4").Value = -0.8236638877349185
$ws.Range("T4").Value = -0.2921931896054355
$ws.Range("U4").Value = 0.01144978231697838
$ws.Range("V4").Value = 0.02951348702483948
$ws.Range("W4").Value = 0.01204206861318278
$ws.Range("X4").Value = -0.006321766063982288
$ws.Range("Y4").Value = -0.09113816209965289
$ws.Range("Z4").Value = 0.08939499002836609
$ws.Range("AA4").Value = 0.01084509490996297
$ws.Range("AB4").Value = 0.02722035516329217
$ws.Range("AC4").Value = 0.02723480434075968
$ws.Range("AD4").Value = 0.01084890508582021
$ws.Range("D5").Value = 6.43
$ws.Range("E5").Value = -73.56999999999999
$ws.Range("F5").Value = -0.5707159890124422
$ws.Range("G5").Value = -0.1174547179070453
$ws.Range("H5").Value = -0.2259271411898982
$ws.Range("I5").Value = 0.126074842896542
$ws.Range("J5").Value = 11.85964098587569
$ws.Range("K5").Value = 11.85964098587569
$ws.Range("L5").Value = 11.85964098587569
$ws.Range("M5").Value = -0.2536395448336352
$ws.Range("N5").Value = 0.1065011544469537
$ws.Range("O5").Value = 11.85964098587569
$ws.Range("P5").Value = 11.85964098587569
$ws.Range("Q5").Value = -0.2259271411898982
$ws.Range("R5").Value = 0.126074842896542
$ws.Range("S5").Value = -0.5707159890124422
$ws.Range("T5").Value = -0.1174547179070453
$ws.Range("U5").Value = 0.009394034799342758
$ws.Range("V5").Value = 0.03329805302032233
$ws.Range("W5").Value = 0.005337187945480913
$ws.Range("X5").Value = -0.002978149341335288
$ws.Range("Y5").Value = -0.04765163095939182
$ws.Range("Z5").Value = 0.07101615487188409
$ws.Range("AA5").Value = 0.01050696060608766
$ws.Range("AB5").Value = 0.02863661498534252
$ws.Range("AC5").Value = 0.02873040712614659
$ws.Range("AD5").Value = 0.01048456758204928
$ws.Range("D6").Value = 2.97
$ws.Range("E6").Value = -57.03
$ws.Range("F6").Value = -0.3955308947251313
$ws.Range("G6").Value = 0.005832105706253421
$ws.Range("H6").Value = -0.1369392588582232
$ws.Range("I6").Value = 0.1884792763089439
$ws.Range("J6").Value = 10.80520911087116
$ws.Range("K6").Value = 10.80520911087116
$ws.Range("L6").Value = 10.80520911087116
$ws.Range("M6").Value = -0.1497395448336352
$ws.Range("N6").Value = 0.1794382413641107
$ws.Range("O6").Value = 10.80520911087116
$ws.Range("P6").Value = 10.80520911087116
$ws.Range("Q6").Value = -0.1369392588582232
$ws.Range("R6").Value = 0.1884792763089439
$ws.Range("S6").Value = -0.3955308947251313
$ws.Range("T6").Value = 0.005832105706253421
$ws.Range("U6").Value = 0.006556296475066355
$ws.Range("V6").Value = 0.03587952163432436
$ws.Range("W6").Value = 0.002718053751004552
$ws.Range("X6").Value = -0.001761795041592241
$ws.Range("Y6").Value = -0.05530007483569668
$ws.Range("Z6").Value = 0.07585499778225262
$ws.Range("AA6").Value = 0.01015411342174358
$ws.Range("AB6").Value = 0.02928350349137697
$ws.Range("AC6").Value = 0.02942538780110549
$ws.Range("AD6").Value = 0.01007672221951429
$ws.Range("D7").Value = -90.70999999999999
$ws.Range("F7").Value = -0.4580241149766267
$ws.Range("G7").Value = -0.2863227593446274
$ws.Range("H7").Value = -0.3146781514944039
$ws.Range("I7").Value = -0.1850753444405359
$ws.Range("J7").Value = 9.092753770228061
$ws.Range("K7").Value = 9.092753770228061
$ws.Range("L7").Value = 9.092753770228061
$ws.Range("M7").Value = -0.05404651494454579
$ws.Range("N7").Value = -0.0009872906030907874
$ws.Range("O7").Value = 9.092753770228061
$ws.Range("P7").Value = 9.092753770228061
$ws.Range("Q7").Value = -0.3146781514944039
$ws.Range("R7").Value = -0.1850753444405359
$ws.Range("S7").Value = -0.4580241149766267
$ws.Range("T7").Value = -0.2863227593446274
$ws.Range("U7").Value = -0.01192397566847313
$ws.Range("V7").Value = 0.04150954924400135
$ws.Range("W7").Value = 0.001499874174637297
$ws.Range("X7").Value = -0.001187609127327498
$ws.Range("Y7").Value = -0.1041505698014761
$ws.Range("Z7").Value = 0.07733289801797667
$ws.Range("AA7").Value = 0.009446238777498353
$ws.Range("AB7").Value = 0.02967549881997178
$ws.Range("AC7").Value = 0.02986536854883751
$ws.Range("AD7").Value = 0.009103367429373951
$ws.Range("D8").Value = 140.6
$ws.Range("E8").Value = -140.6
$ws.Range("F8").Value = -0.4780241149766266
$ws.Range("G8").Value = -0.3293684361214495
$ws.Range("H8").Value = 0.329931085087535
$ws.Range("I8").Value = 0.2413025013616236
$ws.Range("J8").Value = 7.642638181304306
$ws.Range("K8").Value = 7.642638181304306
$ws.Range("L8").Value = 7.642638181304306
$ws.Range("M8").Value = -0.07404651494454578
$ws.Range("N8").Value = -0.04403296737991291
$ws.Range("O8").Value = 7.642638181304306
$ws.Range("P8").Value = 7.642638181304306
$ws.Range("Q8").Value = 0.329931085087535
$ws.Range("R8").Value = 0.2413025013616236
$ws.Range("S8").Value = -0.4780241149766266
$ws.Range("T8").Value = -0.3293684361214495
$ws.Range("U8").Value = -0.07590698325102638
$ws.Range("V8").Value = 0.05086667861553389
$ws.Range("W8").Value = 0.00103434299229522
$ws.Range("X8").Value = -0.0009692518520112886
$ws.Range("Y8").Value = -0.1013799025838043
$ws.Range("Z8").Value = 0.05279486218599103
$ws.Range("AA8").Value = 0.007663097072629941
$ws.Range("AB8").Value = 0.03011821059547105
$ws.Range("AC8").Value = 0.03050868491997199
$ws.Range("AD8").Value = 0.006090355864396359
$ws.Range("D9").Value = 140.6
$ws.Range("E9").Value = -140.6
$ws.Range("F9").Value = -0.4524241149766267
$ws.Range("G9").Value = -0.3121822319566691
$ws.Range("H9").Value = 0.355531085087535
$ws.Range("I9").Value = 0.2584887055264041
$ws.Range("J9").Value = 5.410576205349702
$ws.Range("K9").Value = 5.410576205349702
$ws.Range("L9").Value = 5.410576205349702
$ws.Range("M9").Value = -0.04844651494454579
$ws.Range("N9").Value = -0.02684676321513248
$ws.Range("O9").Value = 5.410576205349702
$ws.Range("P9").Value = 5.410576205349702
$ws.Range("Q9").Value = 0.355531085087535
$ws.Range("R9").Value = 0.2584887055264041
$ws.Range("S9").Value = -0.4524241149766267
$ws.Range("T9").Value = -0.3121822319566691
$ws.Range("U9").Value = -0.09010518232724722
$ws.Range("V9").Value = 0.0471575020577125
$ws.Range("W9").Value = 0.0005937145505192322
$ws.Range("X9").Value = -0.0002968572752596161
$ws.Range("AA9").Value = 0.005417402809721732
$ws.Range("AB9").Value = 0.03050959560618879
$ws.Range("AC9").Value = 0.03100978314244442
$ws.Range("AD9").Value = 0.01007672221951429
